$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.021.36"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "2.462.03"
$ws.Range("E3").Value = "  -1.47%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.557"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("D9").Value = "2.471.96"
$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.157"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.43%  "

$ws.Range("E13").Value = "  -2.89%  "

$ws.Range("D14").Value = "2.902.92"
$ws.Range("E14").Value = "  -2.00%  "

$ws.Range("D15").Value = "57.961.33"
$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.93%  "

$ws.Range("E17").Value = "  -2.63%  "

$ws.Range("D18").Value = "2.465.46"
$ws.Range("E18").Value = "  -1.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "320.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.99%  "

$ws.Range("E25").Value = "  -2.97%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  -3.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.62%  "

$ws.Range("E29").Value = "  -2.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.62%  "

$ws.Range("E31").Value = "  -3.41%  "

$ws.Range("E32").Value = "  -5.64%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.68%  "

$ws.Range("E38").Value = "  -2.20%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.55%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.795"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.78%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.79%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "273.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.588"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "123.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.36%  "

$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0487"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.17%  "

$ws.Range("E49").Value = "  -3.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.78%  "

$ws.Range("D51").Value = "1.720.97"
$ws.Range("E51").Value = "  -2.06%  "
